# Fruta / hortaliza, semanal
# A new weekly record is inserted as the first data row (row 33) of the
# "Hortaliza, Terminal Hortofrutícola Agro Chillán - Arveja Verde" sheet.
# This pushes every existing record down by one row (old row 33 becomes
# row 34, ..., old row 99 becomes row 100), growing the used range from
# A1:R99 to A1:R100.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current first data row (row 33). This shifts
# rows 33:99 down to 34:100, preserving all of their data/formatting.
$ws.Rows(33).Insert()

# Populate the new row 33 with the new weekly observation.
$ws.Cells.Item(33, 1).Value2  = 7
$ws.Cells.Item(33, 2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(33, 3).Value2  = "Ñuble"
$ws.Cells.Item(33, 4).Value2  = 44935
$ws.Cells.Item(33, 5).Value2  = 16
$ws.Cells.Item(33, 6).Value2  = 100112022
$ws.Cells.Item(33, 7).Value2  = "Arveja Verde"
$ws.Cells.Item(33, 8).Value2  = "Sin especificar"
$ws.Cells.Item(33, 9).Value2  = "Primera"
$ws.Cells.Item(33, 10).Value2 = 60
$ws.Cells.Item(33, 11).Value2 = 22000
$ws.Cells.Item(33, 12).Value2 = 23000
$ws.Cells.Item(33, 13).Value2 = 22500
$ws.Cells.Item(33, 14).Value2 = "$/saco 25 kilos"
$ws.Cells.Item(33, 15).Value2 = "Región de Ñuble"
$ws.Cells.Item(33, 16).Value2 = 900
$ws.Cells.Item(33, 17).Value2 = 25
$ws.Cells.Item(33, 18).Value2 = "Hortaliza"
